$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product availability text for row 274 (Naturaplan Bio Artischocke)
$ws.Cells.Item(274, 13).Value = "Naturaplan Bio Artischocke 1 Stück - Online kein Bestand 1.95 Schweizer Franken"

# Update the scrape timestamp (column O) for every data row (2 through 514)
# from the old crawl time to the new crawl time.
for ($r = 2; $r -le 514; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-01-09 16:02:42"
}
